$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.306.08'
$ws.Range("E2").Value = '  +1.29%  '
$ws.Range("D3").Value = '3.162.51'
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.43'
$ws.Range("E5").Value = '  +1.62%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.81'
$ws.Range("E6").Value = '  +1.25%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '3.151.03'
$ws.Range("E8").Value = '  +1.94%  '
$ws.Range("E9").Value = '  +0.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.164'
$ws.Range("E10").Value = '  +3.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.91'
$ws.Range("E11").Value = '  +4.97%  '
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.46'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").Value = '3.680.88'
$ws.Range("E15").Value = '  +1.96%  '
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.27'
$ws.Range("E17").Value = '  +2.12%  '
$ws.Range("D18").Value = '64.098.36'
$ws.Range("D19").Value = '3.154.87'
$ws.Range("E19").Value = '  +1.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.23'
$ws.Range("E20").Value = '  +2.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.40'
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("E22").Value = '  +1.38%  '
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("E24").Value = '  +12.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.18'
$ws.Range("E25").Value = '  +1.57%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  +10.11%  '
$ws.Range("E29").Value = '  +2.16%  '
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.33'
$ws.Range("E31").Value = '  +8.04%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.87'
$ws.Range("E33").Value = '  +4.30%  '
$ws.Range("E34").Value = '  +5.55%  '
$ws.Range("D35").Value = '0.0₃0871'
$ws.Range("E35").Value = '  +2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +3.33%  '
$ws.Range("E37").Value = '  +2.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.29'
$ws.Range("E39").Value = '  -3.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '466.32'
$ws.Range("E40").Value = '  +7.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.37'
$ws.Range("E41").Value = '  +7.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '51.38'
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("E43").Value = '  +9.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0374'
$ws.Range("E44").Value = '  +1.67%  '
$ws.Range("D45").Value = '2.912.45'
$ws.Range("E45").Value = '  +1.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.94'
$ws.Range("E46").Value = '  +11.94%  '
$ws.Range("E47").Value = '  -0.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.17'
$ws.Range("E48").Value = '  +7.51%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  +4.50%  '
